# Update countries & provincias Spain
# - Refresh "Datos actualizados..." timestamp (22:35 -> 23:05)
# - Update case counters for Estados Unidos, Alemania, Peru
# - Update case counters for "Republica de Africa Central" and re-sort it
#   above Nepal / Etiopia (list is kept sorted by "Casos totales" desc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: refresh "datos actualizados" timestamp -----------------------
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 23:05"

# --- Estados Unidos (row 4) -------------------------------------------------
$ws.Cells.Item(4, 2).Value = 1682966   # Casos totales
$ws.Cells.Item(4, 3).Value = 16138     # Nuevos casos
$ws.Cells.Item(4, 4).Value = 451396    # Casos activos
$ws.Cells.Item(4, 5).Value = 1132342   # Recuperados
$ws.Cells.Item(4, 7).Value = 545       # Casos criticos
$ws.Cells.Item(4, 8).Value = 99228     # Muertes

# --- Alemania (row 11) -------------------------------------------------
$ws.Cells.Item(11, 2).Value = 180328   # Casos totales
$ws.Cells.Item(11, 3).Value = 342      # Nuevos casos
$ws.Cells.Item(11, 5).Value = 11657    # Recuperados

# --- Peru (row 15) -------------------------------------------------
$ws.Cells.Item(15, 4).Value = 49795    # Casos activos
$ws.Cells.Item(15, 5).Value = 66708    # Recuperados
$ws.Cells.Item(15, 7).Value = 83       # Casos criticos
$ws.Cells.Item(15, 8).Value = 3456     # Muertes

# --- Republica de Africa Central / Nepal / Etiopia (rows 130-132) ---------
# Previously sorted: Malta(129), Nepal(130), Etiopia(131), Republica de
# Africa Central(132), Jamaica(133). Republica de Africa Central's updated
# "Casos totales" (604) now outranks Nepal (603) and Etiopia (582), so the
# list re-sorts to: Malta, Republica de Africa Central, Nepal, Etiopia,
# Jamaica. Nepal and Etiopia data itself does not change, only its row
# position; Republica de Africa Central also receives new counters.

$ws.Cells.Item(130, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(130, 2).Value = 604
$ws.Cells.Item(130, 3).Value = 52
$ws.Cells.Item(130, 4).Value = 22
$ws.Cells.Item(130, 5).Value = 581
$ws.Cells.Item(130, 6).Value = 0
$ws.Cells.Item(130, 7).Value = 0
$ws.Cells.Item(130, 8).Value = 1

$ws.Cells.Item(131, 1).Value = "Nepal"
$ws.Cells.Item(131, 2).Value = 603
$ws.Cells.Item(131, 3).Value = 19
$ws.Cells.Item(131, 4).Value = 87
$ws.Cells.Item(131, 5).Value = 513
$ws.Cells.Item(131, 6).Value = 0
$ws.Cells.Item(131, 7).Value = 0
$ws.Cells.Item(131, 8).Value = 3

$ws.Cells.Item(132, 1).Value = "Etiopia"
$ws.Cells.Item(132, 2).Value = 582
$ws.Cells.Item(132, 3).Value = 88
$ws.Cells.Item(132, 4).Value = 152
$ws.Cells.Item(132, 5).Value = 425
$ws.Cells.Item(132, 6).Value = 0
$ws.Cells.Item(132, 7).Value = 0
$ws.Cells.Item(132, 8).Value = 5
